$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header + municipality/state name title-casing for Spanish connector words
# (de/del/la/las/el/los/y -> De/Del/La/Las/El/Los/Y)
$changes = @{
    'A1' = 'mx_state'
    'B1' = 'mx_municipality'
    'C1' = 'n_matriculas'
    'D1' = 'pct_matriculas'
    'B5' = 'Pabellón De Arteaga'
    'B6' = 'Rincón De Romos'
    'B7' = 'San José De Gracia'
    'B23' = 'Amatenango De La Frontera'
    'B26' = 'Benemérito De Las Américas'
    'B34' = 'Comitán De Domínguez'
    'B49' = 'Mazapa De Madero'
    'B52' = 'Ocozocoautla De Espinosa'
    'B58' = 'Salto De Agua'
    'B59' = 'San Cristóbal De Las Casas'
    'B84' = 'Hidalgo Del Parral'
    'B90' = 'Valle De Zaragoza'
    'A106' = 'Ciudad De México'
    'B110' = 'Cuajimalpa De Morelos'
    'B122' = 'Coneto De Comonfort'
    'B127' = 'Pánuco De Coronado'
    'B130' = 'San Juan Del Río'
    'A136' = 'Estado De México'
    'B136' = 'Acambay De Ruíz Castañeda'
    'B139' = 'Almoloya De Alquisiras'
    'B140' = 'Almoloya De Juárez'
    'B145' = 'Atizapán De Zaragoza'
    'B151' = 'Chapa De Mota'
    'B154' = 'Coacalco De Berriozábal'
    'B160' = 'Ecatepec De Morelos'
    'B166' = 'Ixtapan De La Sal'
    'B167' = 'Ixtapan Del Oro'
    'B179' = 'Naucalpan De Juárez'
    'B189' = 'San Antonio La Isla'
    'B190' = 'San Felipe Del Progreso'
    'B192' = 'San Simón De Guerrero'
    'B202' = 'Tenango Del Valle'
    'B214' = 'Tlalnepantla De Baz'
    'B220' = 'Valle De Bravo'
    'B221' = 'Valle De Chalco Solidaridad'
    'B222' = 'Villa De Allende'
    'B223' = 'Villa Del Carbón'
    'B234' = 'Apaseo El Alto'
    'B240' = 'Dolores Hidalgo Cuna De La Independencia Nacional'
    'B249' = 'San Francisco Del Rincón'
    'B251' = 'Santa Cruz De Juventino Rosas'
    'B252' = 'Silao De La Victoria'
    'B256' = 'Valle De Santiago'
    'B261' = 'Acapulco De Juárez'
    'B264' = 'Ajuchitlán Del Progreso'
    'B265' = 'Alcozauca De Guerrero'
    'B269' = 'Atenango Del Río'
    'B270' = 'Atlamajalcingo Del Monte'
    'B272' = 'Atoyac De Álvarez'
    'B273' = 'Ayutla De Los Libres'
    'B276' = 'Buenavista De Cuéllar'
    'B277' = 'Chilapa De Álvarez'
    'B278' = 'Chilpancingo De Los Bravo'
    'B279' = 'Coahuayutla De José María Izazaga'
    'B284' = 'Coyuca De Benítez'
    'B285' = 'Coyuca De Catalán'
    'B289' = 'Cuetzala Del Progreso'
    'B290' = 'Cutzamala De Pinzón'
    'B296' = 'Huitzuco De Los Figueroa'
    'B297' = 'Iguala De La Independencia'
    'B299' = 'Ixcateopan De Cuauhtémoc'
    'B300' = 'Zihuatanejo De Azueta'
    'B304' = 'Mártir De Cuilapan'
    'B316' = 'Taxco De Alarcón'
    'B318' = 'Técpan De Galeana'
    'B320' = 'Tepecoacuilco De Trujano'
    'B322' = 'Tixtla De Guerrero'
    'B326' = 'Tlalixtaquilla De Maldonado'
    'B327' = 'Tlapa De Comonfort'
    'B339' = 'Agua Blanca De Iturbide'
    'B344' = 'Atotonilco El Grande'
    'B350' = 'Cuautepec De Hinojosa'
    'B356' = 'Huasca De Ocampo'
    'B360' = 'Huejutla De Reyes'
    'B363' = 'Jacala De Ledezma'
    'B368' = 'Mineral Del Chico'
    'B369' = 'Mineral Del Monte'
    'B370' = 'Mixquiahuala De Juárez'
    'B371' = 'Molango De Escamilla'
    'B373' = 'Nopala De Villagrán'
    'B374' = 'Pachuca De Soto'
    'B376' = 'Progreso De Obregón'
    'B380' = 'Santiago De Anaya'
    'B381' = 'Santiago Tulantepec De Lugo Guerrero'
    'B384' = 'Tenango De Doria'
    'B386' = 'Tepeji Del Río De Ocampo'
    'B387' = 'Tezontepec De Aldama'
    'B393' = 'Tula De Allende'
    'B394' = 'Tulancingo De Bravo'
    'B395' = 'Villa De Tezontepec'
    'B399' = 'Zacualtipán De Ángeles'
    'B406' = 'Autlán De Navarro'
    'B409' = 'Encarnación De Díaz'
    'B414' = 'Lagos De Moreno'
    'B417' = 'Ojuelos De Jalisco'
    'B423' = 'San Juanito De Escobedo'
    'B426' = 'Tamazula De Gordiano'
    'B429' = 'Tizapán El Alto'
    'B430' = 'Tlajomulco De Zúñiga'
    'B434' = 'Unión De San Antonio'
    'B435' = 'Valle De Juárez'
    'B437' = 'Yahualica De González Gallo'
    'B498' = 'Tiquicheo De Nicolás Romero'
    'B520' = 'Coatlán Del Río'
    'B527' = 'Jonacatepec De Leandro Valle'
    'B531' = 'Puente De Ixtla'
    'B537' = 'Tetela Del Volcán'
    'B539' = 'Tlaltizapán De Zapata'
    'B547' = 'Zacualpan De Amilpas'
    'B561' = 'San Nicolás De Los Garza'
    'B564' = 'Acatlán De Pérez Figueroa'
    'B571' = 'Chalcatongo De Hidalgo'
    'B572' = 'Ciénega De Zimatlán'
    'B574' = 'Coicoyán De Las Flores'
    'B575' = 'Constancia Del Rosario'
    'B577' = 'Fresnillo De Trujano'
    'B578' = 'Guadalupe De Ramírez'
    'B579' = 'Heroica Ciudad De Ejutla De Crespo'
    'B580' = 'Heroica Ciudad De Huajuapan De León'
    'B581' = 'Heroica Ciudad De Tlaxiaco'
    'B583' = 'Huautla De Jiménez'
    'B584' = 'Ixtlán De Juárez'
    'B585' = 'Heroica Ciudad De Juchitán De Zaragoza'
    'B590' = 'Mariscala De Juárez'
    'B592' = 'Mazatlán Villa De Flores'
    'B594' = 'Miahuatlán De Porfirio Díaz'
    'B596' = 'Nejapa De Madero'
    'B597' = 'Oaxaca De Juárez'
    'B598' = 'Ocotlán De Morelos'
    'B599' = 'Pinotepa De Don Luis'
    'B600' = 'Putla Villa De Guerrero'
    'B615' = 'San Antonio De La Cal'
    'B635' = 'San José Del Progreso'
    'B639' = 'San Juan Bautista Lo De Soto'
    'B675' = 'San Mateo Del Mar'
    'B684' = 'San Miguel Del Puerto'
    'B685' = 'San Miguel El Grande'
    'B710' = 'San Pedro Y San Pablo Teposcolula'
    'B721' = 'Santa Ana Del Valle'
    'B728' = 'Santa Cruz Tacache De Mina'
    'B732' = 'Santa Inés De Zaragoza'
    'B733' = 'Santa Inés Del Monte'
    'B741' = 'Santa María Del Tule'
    'B792' = 'Tanetze De Zaragoza'
    'B793' = 'Tataltepec De Valdés'
    'B794' = 'Teotitlán De Flores Magón'
    'B795' = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
    'B796' = 'Tlacolula De Matamoros'
    'B798' = 'Tlalixtac De Cabrera'
    'B799' = 'Totontepec Villa De Morelos'
    'B801' = 'Villa De Etla'
    'B802' = 'Villa De Tututepec'
    'B803' = 'Villa De Zaachila'
    'B805' = 'Villa Sola De Vega'
    'B807' = 'Yutanduchi De Guerrero'
    'B808' = 'Zapotitlán Del Río'
    'B811' = 'Zimatlán De Álvarez'
    'B836' = 'Ayotoxco De Guerrero'
    'B842' = 'Chalchicomula De Sesma'
    'B852' = 'Chila De La Sal'
    'B863' = 'Cuapiaxtla De Madero'
    'B867' = 'Cuayuca De Andrade'
    'B868' = 'Cuetzalan Del Progreso'
    'B883' = 'Huehuetlán El Chico'
    'B884' = 'Huehuetlán El Grande'
    'B889' = 'Huitzilan De Serdán'
    'B891' = 'Ixcamilpa De Guerrero'
    'B895' = 'Izúcar De Matamoros'
    'B905' = 'Los Reyes De Juárez'
    'B906' = 'Mazapiltepec De Juárez'
    'B917' = 'Palmar De Bravo'
    'B927' = 'San Diego La Mesa Tochimiltzingo'
    'B943' = 'San Nicolás De Los Ranchos'
    'B947' = 'San Salvador El Seco'
    'B948' = 'San Salvador El Verde'
    'B955' = 'Tecali De Herrera'
    'B963' = 'Tepanco De López'
    'B964' = 'Tepatlaxco De Hidalgo'
    'B970' = 'Tepexi De Rodríguez'
    'B972' = 'Tetela De Ocampo'
    'B973' = 'Teteles De Avila Castillo'
    'B978' = 'Tlacotepec De Benito Juárez'
    'B990' = 'Totoltepec De Guerrero'
    'B992' = 'Tuzamapan De Galeana'
    'B996' = 'Xayacatlán De Bravo'
    'B1002' = 'Xochitlán De Vicente Suárez'
    'B1010' = 'Zapotitlán De Méndez'
    'B1018' = 'Amealco De Bonfil'
    'B1020' = 'Cadereyta De Montes'
    'B1022' = 'Jalpan De Serra'
    'B1023' = 'Landa De Matamoros'
    'B1025' = 'Pinal De Amoles'
    'B1028' = 'San Juan Del Río'
    'B1034' = 'Armadillo De Los Infante'
    'B1035' = 'Axtla De Terrazas'
    'B1037' = 'Ciudad Del Maíz'
    'B1047' = 'Santa María Del Río'
    'B1053' = 'Tanquián De Escobedo'
    'B1078' = 'Nacozari De García'
    'B1090' = 'Jalpa De Méndez'
    'B1109' = 'Acuamanala De Miguel Hidalgo'
    'B1111' = 'Amaxac De Guerrero'
    'B1116' = 'Contla De Juan Cuamatzi'
    'B1122' = 'Ixtacuixtla De Mariano Matamoros'
    'B1126' = 'Mazatecochco De José María Morelos'
    'B1127' = 'Nanacamilpa De Mariano Arista'
    'B1130' = 'Papalotla De Xicohténcatl'
    'B1136' = 'San Pablo Del Monte'
    'B1144' = 'Tepetitla De Lardizábal'
    'B1147' = 'Tetla De La Solidaridad'
    'B1158' = 'Ziltlaltépec De Trinidad Sánchez Santos'
    'B1165' = 'Alto Lucero De Gutiérrez Barrios'
    'B1167' = 'Amatlán De Los Reyes'
    'B1179' = 'Castillo De Teayo'
    'B1193' = 'Cosamaloapan De Carpio'
    'B1208' = 'Hueyapan De Ocampo'
    'B1212' = 'Ixhuatlán De Madero'
    'B1213' = 'Ixhuatlán Del Café'
    'B1222' = 'Juchique De Ferrer'
    'B1226' = 'Lerdo De Tejada'
    'B1229' = 'Martínez De La Torre'
    'B1232' = 'Medellín De Bravo'
    'B1235' = 'Mixtla De Altamirano'
    'B1245' = 'Paso De Ovejas'
    'B1246' = 'Paso Del Macho'
    'B1250' = 'Poza Rica De Hidalgo'
    'B1259' = 'Sayula De Alemán'
    'B1262' = 'Soledad De Doblado'
    'B1289' = 'Vega De Alatorre'
    'B1300' = 'Zontecomatlán De López Y Fuentes'
    'B1301' = 'Zozocolco De Hidalgo'
    'B1312' = 'Jiménez Del Teul'
    'B1315' = 'Nochistlán De Mejía'
    'B1316' = 'Noria De Ángeles'
    'B1321' = 'Villa De Cos'
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}

# Remove the footer/notes rows (1326-1330) entirely; the table now ends at row 1324
$ws.Range("A1326:A1330").EntireRow.Delete()

# Tiny floating-point precision refresh on recomputed percentage cells
$floatFixes = @{
    'D136' = 0.0009626420822688732
    'D274' = 0.0009626420822688732
    'D295' = 0.009663445518160611
    'D545' = 0.0009996667777407533
    'D948' = 0.0009626420822688732
    'D978' = 0.0009626420822688732
    'D1279' = 0.0009996667777407533
}
foreach ($ref in $floatFixes.Keys) {
    $ws.Range($ref).Value = $floatFixes[$ref]
}

Write-Host ("Done. UsedRange: " + $ws.UsedRange.Address())